$wb = $excel.ActiveWorkbook

# --- "posts" sheet: add new TOPIC_URL column (N) with a VLOOKUP formula ---
$ws1 = $wb.Worksheets.Item("posts")
$ws1.Activate()

# Header cell, matching the style used by the other metric headers (I1:M1)
$ws1.Range("M1").Copy($ws1.Range("N1"))
$ws1.Range("N1").Value = "TOPIC_URL"

# Data rows: VLOOKUP the topic name from the "topics" sheet using column C
for ($r = 2; $r -le 11; $r++) {
    $ws1.Cells.Item($r, 14).Formula = "=VLOOKUP(C$r, topics!B`$2:D`$100, 2, FALSE)"
}

# --- view/selection state left behind by the edit ---
$ws2 = $wb.Worksheets.Item("topics")
$ws2.Activate()
$ws2.Rows.Item(12).Select()

$ws1.Activate()
$ws1.Range("M8").Select()
